$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 51, pushing the existing rows 51-121
# down to 52-122 (data, formatting and styles shift along automatically).
$ws.Rows("51:51").Insert()

# Populate the newly inserted row 51 with the new weekly price record.
$ws.Range("A51").Value = 3
$ws.Range("B51").Value = "Femacal de La Calera"
$ws.Range("C51").Value = "Coquimbo"
$ws.Range("D51").Value = 44557
$ws.Range("E51").Value = 5
$ws.Range("F51").Value = 100112026
$ws.Range("G51").Value = "Haba"
$ws.Range("H51").Value = "Sin especificar"
$ws.Range("I51").Value = "Primera"
$ws.Range("J51").Value = 140
$ws.Range("K51").Value = 7000
$ws.Range("L51").Value = 7500
$ws.Range("M51").Value = 7250
$ws.Range("N51").Value = "$/saco 25 kilos"
$ws.Range("O51").Value = "Provincia de Quillota"
$ws.Range("P51").Value = 290
$ws.Range("Q51").Value = 25
$ws.Range("R51").Value = "Hortaliza"
